$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.209.91'
$ws.Range('E2').Value = '  +0.85%  '
$ws.Range('D3').Value = '1.797.50'
$ws.Range('E3').Value = '  +2.03%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = "'338.96"
$ws.Range('E5').Value = '  +0.63%  '
$ws.Range('D6').Value = "'1.001"
$ws.Range('E6').Value = '  +0.17%  '
$ws.Range('D7').Value = "'0.4763"
$ws.Range('E7').Value = '  +26.01%  '
$ws.Range('D8').Value = "'0.3710"
$ws.Range('E8').Value = '  +10.76%  '
$ws.Range('D9').Value = "'45.51"
$ws.Range('E9').Value = '  -0.23%  '
$ws.Range('D10').Value = "'0.07683"
$ws.Range('E10').Value = '  +6.88%  '
$ws.Range('D11').Value = "'1.144"
$ws.Range('D12').Value = "'22.61"
$ws.Range('E12').Value = '  +1.25%  '
$ws.Range('E13').Value = '  +0.08%  '
$ws.Range('D14').Value = "'6.311"
$ws.Range('E14').Value = '  +1.99%  '
$ws.Range('D15').Value = "'7.312"
$ws.Range('E15').Value = '  +1.80%  '
$ws.Range('D16').Value = '1.798.49'
$ws.Range('E16').Value = '  +2.30%  '
$ws.Range('D17').Value = "'0.00001095"
$ws.Range('E17').Value = '  +3.99%  '
$ws.Range('D18').Value = "'0.06725"
$ws.Range('E18').Value = '  +2.29%  '
$ws.Range('D19').Value = "'82.17"
$ws.Range('E19').Value = '  +2.18%  '
$ws.Range('E20').Value = '  +0.06%  '
$ws.Range('E21').Value = '  +2.51%  '
$ws.Range('D22').Value = "'6.409"
$ws.Range('E22').Value = '  +2.03%  '
$ws.Range('D23').Value = '28.208.23'
$ws.Range('E23').Value = '  +0.79%  '
$ws.Range('D24').Value = "'12.01"
$ws.Range('E24').Value = '  +2.70%  '
$ws.Range('D25').Value = "'2.404"
$ws.Range('E25').Value = '  +2.00%  '
$ws.Range('D26').Value = "'20.71"
$ws.Range('E26').Value = '  +4.24%  '
$ws.Range('D27').Value = "'2.406"
$ws.Range('E27').Value = '  +3.06%  '
$ws.Range('D28').Value = "'150.60"
$ws.Range('D29').Value = '2.003.39'
$ws.Range('E29').Value = '  +2.27%  '
$ws.Range('D30').Value = "'133.99"
$ws.Range('E30').Value = '  +1.56%  '
$ws.Range('D31').Value = "'1.271"
$ws.Range('E31').Value = '  +0.44%  '
$ws.Range('D32').Value = "'4.050"
$ws.Range('E32').Value = '  +0.82%  '
$ws.Range('D33').Value = "'0.09652"
$ws.Range('E33').Value = '  +10.04%  '
$ws.Range('D34').Value = "'5.939"
$ws.Range('E34').Value = '  +2.44%  '
$ws.Range('D35').Value = "'0.02372"
$ws.Range('E35').Value = '  +1.65%  '
$ws.Range('D36').Value = "'12.17"
$ws.Range('E36').Value = '  -0.60%  '
$ws.Range('D37').Value = "'0.6685"
$ws.Range('E37').Value = '  +1.56%  '
$ws.Range('B38').Value = 'Algorand'
$ws.Range('C38').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D38').Value = "'0.2190"
$ws.Range('E38').Value = '  +4.03%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').Value = "'0.06310"
$ws.Range('E39').Value = '  +2.03%  '
$ws.Range('B40').Value = 'InternetComputer(DFINITY)'
$ws.Range('C40').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D40').Value = "'5.241"
$ws.Range('E40').Value = '  +1.81%  '
$ws.Range('B41').Value = 'WEMIXTOKEN'
$ws.Range('C41').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D41').Value = "'1.486"
$ws.Range('E41').Value = '  +2.76%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').Value = "'1.222"
$ws.Range('E42').Value = '  +1.05%  '
$ws.Range('D43').Value = "'8.123"
$ws.Range('E43').Value = '  +1.63%  '
$ws.Range('D44').Value = "'14.18"
$ws.Range('E44').Value = '  +2.75%  '
$ws.Range('E45').Value = '  +0.13%  '
$ws.Range('D46').Value = "'0.6158"
$ws.Range('E46').Value = '  +1.92%  '
$ws.Range('D47').Value = "'3.874"
$ws.Range('E47').Value = '  +1.32%  '
$ws.Range('D48').Value = "'129.07"
$ws.Range('E48').Value = '  -0.86%  '
$ws.Range('D49').Value = "'2.047"
$ws.Range('E49').Value = '  +1.85%  '
$ws.Range('D50').Value = "'1.172"
$ws.Range('E50').Value = '  -1.13%  '
$ws.Range('D51').Value = "'0.07095"
$ws.Range('E51').Value = '  -0.98%  '

Write-Output "applied 103 cell updates"
